$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 ("Merge test:") to make room for the new
# "Ignore test:" row, then fill it in with the matching label/command text.
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "Ignore test:"
$ws.Range("B3").Value = "xltablediff.py  --key ID --ignore Color test1old.xlsx test1new.xlsx --out test1ignore.xlsx"
